# DIV2_Tables.xlsx — refresh of the HTotalRevComparison sheet with a newer
# comparison run ("c rev with zip"): PREVIOUS/LATEST columns swap places and
# every partner's figures are updated, which also changes which rows read
# as an increase (green) vs. a decrease (red).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HTotalRevComparison")

# --- Column widths: column B (was the narrower "9" width) now matches the
#     width="10" group shared by columns C and D. ColumnWidth (characters)
#     is offset from the raw sheet width attribute by 5/6, so 9.1666... -> 10.
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# --- Row 1 header labels: PREVIOUS and LATEST swap columns.
$ws.Range("B1").Value = "PREVIOUS"
$ws.Range("C1").Value = "LATEST"

# --- Row 2: BEYOND EXPRESS
$ws.Range("B2").Value = 2641.24
$ws.Range("C2").Value = 1152.02
$ws.Range("D2").Value = -1489.22

# --- Row 3: CIRCLE TRANSIT
$ws.Range("B3").Value = 13278.34
$ws.Range("C3").Value = 13893.26
$ws.Range("D3").Value = 614.92

# --- Row 4: DEAN LOGISTICS SOLUTIONS LLC
#     This partner flips from an increase to a decrease, so the D4 cell's
#     conditional-style swaps from the existing "Increased" (green) look to
#     the existing "Decreased" (red) look. Copy formats only from a cell
#     already carrying the target look so the existing style is reused
#     rather than a new one being minted.
$ws.Range("B4").Value = 2024.62
$ws.Range("C4").Value = 1852.16
$ws.Range("D4").Value = -172.46
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# --- Row 5: DIRECT TRANSIT LLC
$ws.Range("B5").Value = 13231.05
$ws.Range("C5").Value = 14843.35
$ws.Range("D5").Value = 1612.3

# --- Row 6: Q-TELL TRANSIT LLC
#     This partner flips from a decrease to an increase, so D6 picks up the
#     existing green "increase" look (copied from D3, which already has it).
$ws.Range("B6").Value = 4626.1
$ws.Range("C6").Value = 6237.46
$ws.Range("D6").Value = 1611.36
$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
